$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values updated by the automatic electricity price update.
$ws.Range("A2").Value = 45915
$ws.Range("B2").Value = 97.03
$ws.Range("C2").Value = 94.95
$ws.Range("D2").Value = 90
$ws.Range("E2").Value = 89
$ws.Range("F2").Value = 87.81999999999999
$ws.Range("G2").Value = 94.17
$ws.Range("H2").Value = 97.33
$ws.Range("I2").Value = 120.83
$ws.Range("J2").Value = 124
$ws.Range("K2").Value = 84.2
$ws.Range("L2").Value = 35.1
$ws.Range("M2").Value = 7.63
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 4.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 3.52
$ws.Range("R2").Value = 4.96
$ws.Range("S2").Value = 16
$ws.Range("T2").Value = 43.99
$ws.Range("U2").Value = 94.17
$ws.Range("V2").Value = 110
$ws.Range("W2").Value = 110.96
$ws.Range("X2").Value = 98.5
$ws.Range("Y2").Value = 85.97
$ws.Range("Z2").Value = 66.88
$ws.Range("AB2").Value = 101.36
$ws.Range("AD2").Value = 110.48
$ws.Range("AE2").Value = "6h-8h"
$ws.Range("AF2").Value = 109.08
$ws.Range("AG2").Value = "10h-18h"
